$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = -0.0005547603745488239
$ws.Range("C2").Value = -0.01415192358968653
$ws.Range("D2").Value = -0.000381023133439129

$ws.Range("B3").Value = 0.004348273989535301
$ws.Range("C3").Value = 0.08422249561067474
$ws.Range("D3").Value = 0.00608758394905351

$ws.Range("B4").Value = 0.002293105418340602
$ws.Range("C4").Value = 8.26788768566712
$ws.Range("D4").Value = 0.02035590906120355

$ws.Range("B5").Value = -0.0001845260572519436
$ws.Range("C5").Value = -0.002499387409475418
$ws.Range("D5").Value = -0.00006120791099673006

$ws.Range("B6").Value = -0.0005381055849119631
$ws.Range("C6").Value = -0.01372706031328885
$ws.Range("D6").Value = -0.0003695842123744297

$ws.Range("B7").Value = -0.05678528333390886
$ws.Range("C7").Value = -142.8662955041
$ws.Range("D7").Value = -0.5040832638198083

$ws.Range("B8").Value = -0.1598975545786061
$ws.Range("C8").Value = -9.133492215009028
$ws.Range("D8").Value = -0.09902364618739057

$ws.Range("B9").Value = -0.0002819448126429869
$ws.Range("C9").Value = -0.008701059050579829
$ws.Range("D9").Value = -0.01609815991332653

$ws.Range("B10").Value = -0.0779959105711896
$ws.Range("C10").Value = -1.276521816942477
$ws.Range("D10").Value = -0.08850654657362611

$ws.Range("B11").Value = -0.00006420625365866073
$ws.Range("D11").Value = -0.001123360939317308
